$d = $word.ActiveDocument

function Split-Run {
    param(
        [string]$NeedleText
    )
    # Locate the (now up-to-date) text and force Word to materialize it as
    # its own run by toggling a character property on then off. This mirrors
    # how Word splits a run when a sub-span of text gets distinct formatting
    # applied and then removed again.
    $r = $d.Content.Duplicate()
    $r.Find.Execute($NeedleText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Bold = 1
    $r.Bold = 0
}

# --- Edit 1: Objectives bullet -------------------------------------------
# "Build an AI-powered application using YOLO for object detection"
# -> "Build an AI-powered application using YOLOv8 for object detection"
# split into 3 runs around "YOLOv8"
$rng1 = $d.Content.Duplicate()
$rng1.Find.Execute("Build an AI-powered application using YOLO for object detection", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Build an AI-powered application using YOLOv8 for object detection", 2)

Split-Run "YOLOv8 for object detection"

# --- Edit 2: Technologies paragraph --------------------------------------
# "YOLO" (standalone run) -> "YOLOv8"
$rng2 = $d.Content.Duplicate()
$rng2.Find.Execute("YOLO which will be built on PyTorch", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(1)
$rng2.MoveEnd(1, 4)
$rng2.Text = "YOLOv8"

# --- Edit 3: Steps table cell ---------------------------------------------
# "Train YOLO on annotated datasets using college resources"
# -> "Train YOLOv8 on annotated datasets using college resources"
# split into 3 runs around "YOLOv8"
$rng3 = $d.Content.Duplicate()
$rng3.Find.Execute("Train YOLO on annotated datasets using college resources", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Train YOLOv8 on annotated datasets using college resources", 2)

Split-Run "Train "
Split-Run "YOLOv8 on annotated datasets using college resources"

# --- Edit 4: Budget bullet --------------------------------------------------
# "No cost due to open-source and free tools like PyTorch, YOLO, Groq"
# -> "No cost due to open-source and free tools like PyTorch, YOLOv8, Groq"
# split into 3 runs around "YOLOv8"
$rng4 = $d.Content.Duplicate()
$rng4.Find.Execute("No cost due to open-source and free tools like PyTorch, YOLO, Groq", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "No cost due to open-source and free tools like PyTorch, YOLOv8, Groq", 2)

Split-Run "No cost due to open-source and free tools like PyTorch, "
Split-Run "YOLOv8, Groq"
